{"js": "// The document has four \"Uke 3 <weekday> <date>\" sections, each followed\n// by an empty paragraph and a \"[placeholder]\" paragraph. The edit targets\n// only the FIRST such placeholder \u2014 the one right after the\n// \"Uke 3 Tirsdag 17.1.2023\" heading \u2014 and replaces it with the actual\n// documentation text for that day.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the heading paragraph for \"Uke 3 Tirsdag 17.1.2023\" and then the\n// very next paragraph whose text is exactly \"[placeholder]\" \u2014 that is the\n// paragraph this diff rewrites.\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Tirsdag 17.1.2023\") !== -1) {\n    for (let j = i + 1; j < paragraphs.items.length; j++) {\n      const candidate = paragraphs.items[j].text.trim();\n      if (candidate === \"[placeholder]\") {\n        targetParagraph = paragraphs.items[j];\n        break;\n      }\n    }\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the [placeholder] paragraph for Tirsdag 17.1.2023\");\n}\n\nconst newText =\n  \"Viet-Uy, Besnik og Ilyas kom f\u00f8rst inn i timen og rett etter kom Mathias. \" +\n  \"Viet-Uy begynte med \u00e5 sette opp en router. \" +\n  \"Ilyas jobbet mer med nettside skisse, \";\n\ntargetParagraph.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document has four \"Uke 3 <weekday> <date>\" sections, each followed\n# by an empty paragraph and a \"[placeholder]\" paragraph. This edit targets\n# only the FIRST such placeholder \u2014 the one right after the\n# \"Uke 3 Tirsdag 17.1.2023\" heading \u2014 and replaces it with the actual\n# documentation text for that day.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*Tirsdag 17.1.2023*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the 'Uke 3 Tirsdag 17.1.2023' heading paragraph\"\n}\n\n# Walk forward from the heading to the next paragraph whose text is the\n# literal placeholder \"[placeholder]\".\n$placeholderIndex = -1\nfor ($j = $targetIndex + 1; $j -le $d.Paragraphs.Count; $j++) {\n    $candidate = $d.Paragraphs.Item($j).Range.Text.Trim()\n    if ($candidate -eq \"[placeholder]\") {\n        $placeholderIndex = $j\n        break\n    }\n}\n\nif ($placeholderIndex -eq -1) {\n    throw \"Could not locate the [placeholder] paragraph for Tirsdag 17.1.2023\"\n}\n\n$newText = \"Viet-Uy, Besnik og Ilyas kom f\u00f8rst inn i timen og rett etter kom Mathias. \" + `\n    \"Viet-Uy begynte med \u00e5 sette opp en router. \" + `\n    \"Ilyas jobbet mer med nettside skisse, \"\n\n$range = $d.Paragraphs.Item($placeholderIndex).Range\n$find = $range.Find\n$find.Text = \"[placeholder]\"\n$find.Replacement.Text = $newText\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
